$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G4").Value = "Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
Write-Host "Cell G4 now: $($ws.Range('G4').Value)"
